# Natmi following Dr Hou advice
# Update LR-pair stats for Anxa1-Dysf: sending/target cluster counts go from 1 to 3 replicates,
# recompute derived metrics accordingly, and add the missing sCs sending-cluster rows (14-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Anxa1"
$ws.Cells.Item(2, 3).Value = "Dysf"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"50.24752933333333"
$ws.Cells.Item(2, 8).Value = [double]"150.742588"
$ws.Cells.Item(2, 9).Value = [double]"0.1294604271951564"
$ws.Cells.Item(2, 10).Value = [double]"0.1294604271951564"
$ws.Cells.Item(2, 11).Value = [double]"3"
$ws.Cells.Item(2, 12).Value = [double]"1"
$ws.Cells.Item(2, 13).Value = [double]"27.29291933333333"
$ws.Cells.Item(2, 14).Value = [double]"81.878758"
$ws.Cells.Item(2, 15).Value = [double]"0.8471655416169349"
$ws.Cells.Item(2, 16).Value = [double]"0.8471655416169349"
$ws.Cells.Item(2, 17).Value = [double]"1371.401764793967"
$ws.Cells.Item(2, 18).Value = [double]"12342.6158831457"
$ws.Cells.Item(2, 19).Value = [double]"0.1096744129227444"
$ws.Cells.Item(2, 20).Value = [double]"0.1096744129227444"

# row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Anxa1"
$ws.Cells.Item(3, 3).Value = "Dysf"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"50.24752933333333"
$ws.Cells.Item(3, 8).Value = [double]"150.742588"
$ws.Cells.Item(3, 9).Value = [double]"0.1294604271951564"
$ws.Cells.Item(3, 10).Value = [double]"0.1294604271951564"
$ws.Cells.Item(3, 11).Value = [double]"3"
$ws.Cells.Item(3, 12).Value = [double]"1"
$ws.Cells.Item(3, 13).Value = [double]"1.664698666666667"
$ws.Cells.Item(3, 14).Value = [double]"4.994096"
$ws.Cells.Item(3, 15).Value = [double]"0.0516718395108896"
$ws.Cells.Item(3, 16).Value = [double]"0.0516718395108896"
$ws.Cells.Item(3, 17).Value = [double]"83.64699508449421"
$ws.Cells.Item(3, 18).Value = [double]"752.8229557604479"
$ws.Cells.Item(3, 19).Value = [double]"0.006689458417039326"
$ws.Cells.Item(3, 20).Value = [double]"0.006689458417039328"

# row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Anxa1"
$ws.Cells.Item(4, 3).Value = "Dysf"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"50.24752933333333"
$ws.Cells.Item(4, 8).Value = [double]"150.742588"
$ws.Cells.Item(4, 9).Value = [double]"0.1294604271951564"
$ws.Cells.Item(4, 10).Value = [double]"0.1294604271951564"
$ws.Cells.Item(4, 11).Value = [double]"2"
$ws.Cells.Item(4, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4, 13).Value = [double]"0.01065033333333333"
$ws.Cells.Item(4, 14).Value = [double]"0.031951"
$ws.Cells.Item(4, 15).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(4, 16).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(4, 17).Value = [double]"0.5351529365764444"
$ws.Cells.Item(4, 18).Value = [double]"4.816376429188"
$ws.Cells.Item(4, 19).Value = [double]"4.279751247930026E-05"
$ws.Cells.Item(4, 20).Value = [double]"4.279751247930027E-05"

# row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Anxa1"
$ws.Cells.Item(5, 3).Value = "Dysf"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"50.24752933333333"
$ws.Cells.Item(5, 8).Value = [double]"150.742588"
$ws.Cells.Item(5, 9).Value = [double]"0.1294604271951564"
$ws.Cells.Item(5, 10).Value = [double]"0.1294604271951564"
$ws.Cells.Item(5, 11).Value = [double]"3"
$ws.Cells.Item(5, 12).Value = [double]"1"
$ws.Cells.Item(5, 13).Value = [double]"3.248480333333333"
$ws.Cells.Item(5, 14).Value = [double]"9.745441"
$ws.Cells.Item(5, 15).Value = [double]"0.1008320351300503"
$ws.Cells.Item(5, 16).Value = [double]"0.1008320351300502"
$ws.Cells.Item(5, 17).Value = [double]"163.2281108379231"
$ws.Cells.Item(5, 18).Value = [double]"1469.052997541308"
$ws.Cells.Item(5, 19).Value = [double]"0.01305375834289332"
$ws.Cells.Item(5, 20).Value = [double]"0.01305375834289332"

# row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Anxa1"
$ws.Cells.Item(6, 3).Value = "Dysf"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = [double]"3"
$ws.Cells.Item(6, 6).Value = [double]"1"
$ws.Cells.Item(6, 7).Value = [double]"151.42276"
$ws.Cells.Item(6, 8).Value = [double]"454.26828"
$ws.Cells.Item(6, 9).Value = [double]"0.3901337131747328"
$ws.Cells.Item(6, 10).Value = [double]"0.3901337131747328"
$ws.Cells.Item(6, 11).Value = [double]"3"
$ws.Cells.Item(6, 12).Value = [double]"1"
$ws.Cells.Item(6, 13).Value = [double]"27.29291933333333"
$ws.Cells.Item(6, 14).Value = [double]"81.878758"
$ws.Cells.Item(6, 15).Value = [double]"0.8471655416169349"
$ws.Cells.Item(6, 16).Value = [double]"0.8471655416169349"
$ws.Cells.Item(6, 17).Value = [double]"4132.769173910694"
$ws.Cells.Item(6, 18).Value = [double]"37194.92256519624"
$ws.Cells.Item(6, 19).Value = [double]"0.3305078384246984"
$ws.Cells.Item(6, 20).Value = [double]"0.3305078384246984"

# row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Anxa1"
$ws.Cells.Item(7, 3).Value = "Dysf"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = [double]"3"
$ws.Cells.Item(7, 6).Value = [double]"1"
$ws.Cells.Item(7, 7).Value = [double]"151.42276"
$ws.Cells.Item(7, 8).Value = [double]"454.26828"
$ws.Cells.Item(7, 9).Value = [double]"0.3901337131747328"
$ws.Cells.Item(7, 10).Value = [double]"0.3901337131747328"
$ws.Cells.Item(7, 11).Value = [double]"3"
$ws.Cells.Item(7, 12).Value = [double]"1"
$ws.Cells.Item(7, 13).Value = [double]"1.664698666666667"
$ws.Cells.Item(7, 14).Value = [double]"4.994096"
$ws.Cells.Item(7, 15).Value = [double]"0.0516718395108896"
$ws.Cells.Item(7, 16).Value = [double]"0.0516718395108896"
$ws.Cells.Item(7, 17).Value = [double]"252.0732666749867"
$ws.Cells.Item(7, 18).Value = [double]"2268.65940007488"
$ws.Cells.Item(7, 19).Value = [double]"0.02015892661495223"
$ws.Cells.Item(7, 20).Value = [double]"0.02015892661495223"

# row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Anxa1"
$ws.Cells.Item(8, 3).Value = "Dysf"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = [double]"3"
$ws.Cells.Item(8, 6).Value = [double]"1"
$ws.Cells.Item(8, 7).Value = [double]"151.42276"
$ws.Cells.Item(8, 8).Value = [double]"454.26828"
$ws.Cells.Item(8, 9).Value = [double]"0.3901337131747328"
$ws.Cells.Item(8, 10).Value = [double]"0.3901337131747328"
$ws.Cells.Item(8, 11).Value = [double]"2"
$ws.Cells.Item(8, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(8, 13).Value = [double]"0.01065033333333333"
$ws.Cells.Item(8, 14).Value = [double]"0.031951"
$ws.Cells.Item(8, 15).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(8, 16).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(8, 17).Value = [double]"1.612702868253333"
$ws.Cells.Item(8, 18).Value = [double]"14.51432581428"
$ws.Cells.Item(8, 19).Value = [double]"0.000128971862830498"
$ws.Cells.Item(8, 20).Value = [double]"0.000128971862830498"

# row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Anxa1"
$ws.Cells.Item(9, 3).Value = "Dysf"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = [double]"3"
$ws.Cells.Item(9, 6).Value = [double]"1"
$ws.Cells.Item(9, 7).Value = [double]"151.42276"
$ws.Cells.Item(9, 8).Value = [double]"454.26828"
$ws.Cells.Item(9, 9).Value = [double]"0.3901337131747328"
$ws.Cells.Item(9, 10).Value = [double]"0.3901337131747328"
$ws.Cells.Item(9, 11).Value = [double]"3"
$ws.Cells.Item(9, 12).Value = [double]"1"
$ws.Cells.Item(9, 13).Value = [double]"3.248480333333333"
$ws.Cells.Item(9, 14).Value = [double]"9.745441"
$ws.Cells.Item(9, 15).Value = [double]"0.1008320351300503"
$ws.Cells.Item(9, 16).Value = [double]"0.1008320351300502"
$ws.Cells.Item(9, 17).Value = [double]"491.8938578790534"
$ws.Cells.Item(9, 18).Value = [double]"4427.04472091148"
$ws.Cells.Item(9, 19).Value = [double]"0.03933797627225161"
$ws.Cells.Item(9, 20).Value = [double]"0.0393379762722516"

# row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Anxa1"
$ws.Cells.Item(10, 3).Value = "Dysf"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = [double]"3"
$ws.Cells.Item(10, 6).Value = [double]"1"
$ws.Cells.Item(10, 7).Value = [double]"85.147481"
$ws.Cells.Item(10, 8).Value = [double]"255.442443"
$ws.Cells.Item(10, 9).Value = [double]"0.2193785328573129"
$ws.Cells.Item(10, 10).Value = [double]"0.2193785328573129"
$ws.Cells.Item(10, 11).Value = [double]"3"
$ws.Cells.Item(10, 12).Value = [double]"1"
$ws.Cells.Item(10, 13).Value = [double]"27.29291933333333"
$ws.Cells.Item(10, 14).Value = [double]"81.878758"
$ws.Cells.Item(10, 15).Value = [double]"0.8471655416169349"
$ws.Cells.Item(10, 16).Value = [double]"0.8471655416169349"
$ws.Cells.Item(10, 17).Value = [double]"2323.923330369533"
$ws.Cells.Item(10, 18).Value = [double]"20915.30997332579"
$ws.Cells.Item(10, 19).Value = [double]"0.185849933607194"
$ws.Cells.Item(10, 20).Value = [double]"0.185849933607194"

# row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Anxa1"
$ws.Cells.Item(11, 3).Value = "Dysf"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = [double]"3"
$ws.Cells.Item(11, 6).Value = [double]"1"
$ws.Cells.Item(11, 7).Value = [double]"85.147481"
$ws.Cells.Item(11, 8).Value = [double]"255.442443"
$ws.Cells.Item(11, 9).Value = [double]"0.2193785328573129"
$ws.Cells.Item(11, 10).Value = [double]"0.2193785328573129"
$ws.Cells.Item(11, 11).Value = [double]"3"
$ws.Cells.Item(11, 12).Value = [double]"1"
$ws.Cells.Item(11, 13).Value = [double]"1.664698666666667"
$ws.Cells.Item(11, 14).Value = [double]"4.994096"
$ws.Cells.Item(11, 15).Value = [double]"0.0516718395108896"
$ws.Cells.Item(11, 16).Value = [double]"0.0516718395108896"
$ws.Cells.Item(11, 17).Value = [double]"141.7448980907253"
$ws.Cells.Item(11, 18).Value = [double]"1275.704082816528"
$ws.Cells.Item(11, 19).Value = [double]"0.01133569234193749"
$ws.Cells.Item(11, 20).Value = [double]"0.01133569234193749"

# row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Anxa1"
$ws.Cells.Item(12, 3).Value = "Dysf"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = [double]"3"
$ws.Cells.Item(12, 6).Value = [double]"1"
$ws.Cells.Item(12, 7).Value = [double]"85.147481"
$ws.Cells.Item(12, 8).Value = [double]"255.442443"
$ws.Cells.Item(12, 9).Value = [double]"0.2193785328573129"
$ws.Cells.Item(12, 10).Value = [double]"0.2193785328573129"
$ws.Cells.Item(12, 11).Value = [double]"2"
$ws.Cells.Item(12, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(12, 13).Value = [double]"0.01065033333333333"
$ws.Cells.Item(12, 14).Value = [double]"0.031951"
$ws.Cells.Item(12, 15).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(12, 16).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(12, 17).Value = [double]"0.9068490551436666"
$ws.Cells.Item(12, 18).Value = [double]"8.161641496293001"
$ws.Cells.Item(12, 19).Value = [double]"7.252297633390404E-05"
$ws.Cells.Item(12, 20).Value = [double]"7.252297633390404E-05"

# row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Anxa1"
$ws.Cells.Item(13, 3).Value = "Dysf"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = [double]"3"
$ws.Cells.Item(13, 6).Value = [double]"1"
$ws.Cells.Item(13, 7).Value = [double]"85.147481"
$ws.Cells.Item(13, 8).Value = [double]"255.442443"
$ws.Cells.Item(13, 9).Value = [double]"0.2193785328573129"
$ws.Cells.Item(13, 10).Value = [double]"0.2193785328573129"
$ws.Cells.Item(13, 11).Value = [double]"3"
$ws.Cells.Item(13, 12).Value = [double]"1"
$ws.Cells.Item(13, 13).Value = [double]"3.248480333333333"
$ws.Cells.Item(13, 14).Value = [double]"9.745441"
$ws.Cells.Item(13, 15).Value = [double]"0.1008320351300503"
$ws.Cells.Item(13, 16).Value = [double]"0.1008320351300502"
$ws.Cells.Item(13, 17).Value = [double]"276.5999174613736"
$ws.Cells.Item(13, 18).Value = [double]"2489.399257152363"
$ws.Cells.Item(13, 19).Value = [double]"0.02212038393184746"
$ws.Cells.Item(13, 20).Value = [double]"0.02212038393184746"

# row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Anxa1"
$ws.Cells.Item(14, 3).Value = "Dysf"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = [double]"3"
$ws.Cells.Item(14, 6).Value = [double]"1"
$ws.Cells.Item(14, 7).Value = [double]"101.3126446666667"
$ws.Cells.Item(14, 8).Value = [double]"303.937934"
$ws.Cells.Item(14, 9).Value = [double]"0.2610273267727979"
$ws.Cells.Item(14, 10).Value = [double]"0.2610273267727979"
$ws.Cells.Item(14, 11).Value = [double]"3"
$ws.Cells.Item(14, 12).Value = [double]"1"
$ws.Cells.Item(14, 13).Value = [double]"27.29291933333333"
$ws.Cells.Item(14, 14).Value = [double]"81.878758"
$ws.Cells.Item(14, 15).Value = [double]"0.8471655416169349"
$ws.Cells.Item(14, 16).Value = [double]"0.8471655416169349"
$ws.Cells.Item(14, 17).Value = [double]"2765.117838333997"
$ws.Cells.Item(14, 18).Value = [double]"24886.06054500598"
$ws.Cells.Item(14, 19).Value = [double]"0.221133356662298"
$ws.Cells.Item(14, 20).Value = [double]"0.221133356662298"

# row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Anxa1"
$ws.Cells.Item(15, 3).Value = "Dysf"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = [double]"3"
$ws.Cells.Item(15, 6).Value = [double]"1"
$ws.Cells.Item(15, 7).Value = [double]"101.3126446666667"
$ws.Cells.Item(15, 8).Value = [double]"303.937934"
$ws.Cells.Item(15, 9).Value = [double]"0.2610273267727979"
$ws.Cells.Item(15, 10).Value = [double]"0.2610273267727979"
$ws.Cells.Item(15, 11).Value = [double]"3"
$ws.Cells.Item(15, 12).Value = [double]"1"
$ws.Cells.Item(15, 13).Value = [double]"1.664698666666667"
$ws.Cells.Item(15, 14).Value = [double]"4.994096"
$ws.Cells.Item(15, 15).Value = [double]"0.0516718395108896"
$ws.Cells.Item(15, 16).Value = [double]"0.0516718395108896"
$ws.Cells.Item(15, 17).Value = [double]"168.6550244930738"
$ws.Cells.Item(15, 18).Value = [double]"1517.895220437664"
$ws.Cells.Item(15, 19).Value = [double]"0.01348776213696055"
$ws.Cells.Item(15, 20).Value = [double]"0.01348776213696055"

# row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Anxa1"
$ws.Cells.Item(16, 3).Value = "Dysf"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = [double]"3"
$ws.Cells.Item(16, 6).Value = [double]"1"
$ws.Cells.Item(16, 7).Value = [double]"101.3126446666667"
$ws.Cells.Item(16, 8).Value = [double]"303.937934"
$ws.Cells.Item(16, 9).Value = [double]"0.2610273267727979"
$ws.Cells.Item(16, 10).Value = [double]"0.2610273267727979"
$ws.Cells.Item(16, 11).Value = [double]"2"
$ws.Cells.Item(16, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(16, 13).Value = [double]"0.01065033333333333"
$ws.Cells.Item(16, 14).Value = [double]"0.031951"
$ws.Cells.Item(16, 15).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(16, 16).Value = [double]"0.0003305837421251881"
$ws.Cells.Item(16, 17).Value = [double]"1.079013436581556"
$ws.Cells.Item(16, 18).Value = [double]"9.711120929234001"
$ws.Cells.Item(16, 19).Value = [double]"8.629139048148585E-05"
$ws.Cells.Item(16, 20).Value = [double]"8.629139048148585E-05"

# row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Anxa1"
$ws.Cells.Item(17, 3).Value = "Dysf"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = [double]"3"
$ws.Cells.Item(17, 6).Value = [double]"1"
$ws.Cells.Item(17, 7).Value = [double]"101.3126446666667"
$ws.Cells.Item(17, 8).Value = [double]"303.937934"
$ws.Cells.Item(17, 9).Value = [double]"0.2610273267727979"
$ws.Cells.Item(17, 10).Value = [double]"0.2610273267727979"
$ws.Cells.Item(17, 11).Value = [double]"3"
$ws.Cells.Item(17, 12).Value = [double]"1"
$ws.Cells.Item(17, 13).Value = [double]"3.248480333333333"
$ws.Cells.Item(17, 14).Value = [double]"9.745441"
$ws.Cells.Item(17, 15).Value = [double]"0.1008320351300503"
$ws.Cells.Item(17, 16).Value = [double]"0.1008320351300502"
$ws.Cells.Item(17, 17).Value = [double]"329.112133717655"
$ws.Cells.Item(17, 18).Value = [double]"2962.009203458894"
$ws.Cells.Item(17, 19).Value = [double]"0.02631991658305787"
$ws.Cells.Item(17, 20).Value = [double]"0.02631991658305787"

